# Apply updated loading_percent values for rows 2-25 (case with 380 kV)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$BValues = @(12.87980402204886, 12.62173965962794, 12.4631488333777, 12.39857926064791, 12.38786373680971, 12.46227766952088, 12.79090681948312, 13.43031670122487, 13.89197197271814, 14.09923321694249, 14.17724546116809, 14.16046635986646, 14.10566124222259, 14.07202768263148, 13.87836489681737, 13.75879828222264, 13.68977243708521, 13.66636024829817, 13.77155324998726, 14.12177227091859, 14.34787032582557, 14.22747737060235, 13.7657876120564, 13.25841356802987)
$CValues = @(10.77420178013385, 10.80531747781169, 10.82537905994657, 10.83379558238694, 10.83520773388333, 10.82549159022508, 10.78473242990555, 10.71235749806065, 10.66373921947658, 10.64260022969361, 10.63473525313751, 10.63642290542861, 10.64195037347762, 10.64535430593128, 10.66514031674842, 10.67752831267262, 10.68474562784416, 10.68720512069315, 10.67620006575076, 10.64032303076421, 10.61769044268825, 10.62969552378421, 10.67680026929711, 10.73113328978827)
$DValues = @(5.87748155105499, 5.840504097921566, 5.817303545317124, 5.807727252075199, 5.80612983990519, 5.817174885722271, 5.864836454092131, 5.954238347038837, 6.017282477346654, 6.045355836538032, 6.055896514054294, 6.053630444664875, 6.046224847018192, 6.041676883682077, 6.015435371759498, 5.999179769960405, 5.989773139613654, 5.986578546054065, 6.000916110852521, 6.048402519308291, 6.07891113869288, 6.062677264867325, 6.000131301114235, 5.930504593633605)
$FValues = @(28.54488774583529, 28.57749054321913, 28.60461493480994, 28.61745223636242, 28.61969150031894, 28.60478084480571, 28.55465272697023, 28.51283881538927, 28.51664873238188, 28.52588270614577, 28.53045661970259, 28.52942366857564, 28.5262374272147, 28.52442598936014, 28.51619610746087, 28.5130677772765, 28.51197459941045, 28.51172577861891, 28.51332771618344, 28.52714408415459, 28.54245137047489, 28.53370789420734, 28.51320800059283, 28.51809011513916)
$GValues = @(3.661674947055762, 3.663616051678476, 3.664871243099651, 3.665398722756004, 3.665487276920488, 3.664878292110081, 3.662331122910706, 3.657836454202055, 3.654836017808041, 3.653535890407939, 3.65305283039787, 3.653156454458791, 3.653495963248314, 3.653705128051722, 3.654922283900399, 3.655685529926691, 3.656130629569765, 3.656282381877, 3.655603650018509, 3.653395990066913, 3.652007168284715, 3.652743481694585, 3.655640648299462, 3.658999151740743)
$KValues = @(8.830765766264237, 8.648728662671546, 8.536390530160377, 8.490532709160039, 8.482915159887824, 8.535772310901796, 8.768154412798285, 9.216677331410867, 9.538412214660342, 9.682418193921171, 9.736559862421013, 9.724917633908772, 9.686880568693081, 9.663529459233096, 9.528949224043044, 9.445747977990596, 9.397674258026782, 9.381361421272585, 9.454627900065397, 9.698063965821595, 9.854867006659086, 9.77140468870847, 9.450614031527408, 9.096472149713309)
$LValues = @(10.77683158768227, 10.75171704551102, 10.73846106199048, 10.73360733530506, 10.7328345922585, 10.73839337849705, 10.76772536437605, 10.842225541638, 10.90702274087913, 10.93861397500853, 10.95087405853357, 10.9482205282022, 10.93961670489286, 10.93438510222898, 10.90500015042794, 10.88751026677025, 10.8776500325597, 10.87434599195319, 10.88935149760982, 10.94213584852801, 10.97836209545891, 10.95887169101363, 10.8885184694382, 10.82028140263798)
$OValues = @(25.58691595381615, 25.65714963746957, 25.7055642473884, 25.72662126158459, 25.73019786700342, 25.70584285775335, 25.61003279919669, 25.46424808165922, 25.38295780830567, 25.35161218369758, 25.34055462194769, 25.34289990328695, 25.35068617884212, 25.35556135119966, 25.38511987845148, 25.40469765898665, 25.41648843386942, 25.42057157084643, 25.40255868514686, 25.3483770974759, 25.31770215072643, 25.33363993749774, 25.40352404809334, 25.49916322995003)

$startRow = 2
$columns = @("B", "C", "D", "F", "G", "K", "L", "O")
$columnData = @{
    "B" = $BValues
    "C" = $CValues
    "D" = $DValues
    "F" = $FValues
    "G" = $GValues
    "K" = $KValues
    "L" = $LValues
    "O" = $OValues
}

foreach ($col in $columns) {
    $vals = $columnData[$col]
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $row = $startRow + $i
        $ws.Range("$col$row").Value = $vals[$i]
    }
}

Write-Host "Updated loading_percent values for rows $startRow-$($startRow + $BValues.Count - 1)"